# Updates the "Pais" sheet: refreshes the covid stats table (columns B-H)
# and re-sorts it by "Casos totales" (column B, descending), which also
# moves some country names (column A) to new rows. Also bumps the
# "Datos actualizados" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 03:20"

$data = New-Object 'object[,]' 206,8
$data[0,0] = "Estados Unidos"
$data[0,1] = 215071
$data[0,2] = 68
$data[0,3] = 8878
$data[0,4] = 201087
$data[0,5] = 5005
$data[0,6] = 4
$data[0,7] = 5106
$data[1,0] = "Italia"
$data[1,1] = 110574
$data[1,2] = 0
$data[1,3] = 16847
$data[1,4] = 80572
$data[1,5] = 4035
$data[1,6] = 0
$data[1,7] = 13155
$data[2,0] = "España"
$data[2,1] = 104118
$data[2,2] = 0
$data[2,3] = 22647
$data[2,4] = 72084
$data[2,5] = 5872
$data[2,6] = 0
$data[2,7] = 9387
$data[3,0] = "China"
$data[3,1] = 81554
$data[3,2] = 0
$data[3,3] = 76238
$data[3,4] = 2004
$data[3,5] = 466
$data[3,6] = 0
$data[3,7] = 3312
$data[4,0] = "Alemania"
$data[4,1] = 77981
$data[4,2] = 0
$data[4,3] = 18700
$data[4,4] = 58350
$data[4,5] = 3408
$data[4,6] = 0
$data[4,7] = 931
$data[5,0] = "Francia"
$data[5,1] = 56989
$data[5,2] = 0
$data[5,3] = 10935
$data[5,4] = 42022
$data[5,5] = 6017
$data[5,6] = 0
$data[5,7] = 4032
$data[6,0] = "Iran"
$data[6,1] = 47593
$data[6,2] = 0
$data[6,3] = 15473
$data[6,4] = 29084
$data[6,5] = 3871
$data[6,6] = 0
$data[6,7] = 3036
$data[7,0] = "Reino Unido"
$data[7,1] = 29474
$data[7,2] = 0
$data[7,3] = 135
$data[7,4] = 26987
$data[7,5] = 163
$data[7,6] = 0
$data[7,7] = 2352
$data[8,0] = "Suiza"
$data[8,1] = 17768
$data[8,2] = 0
$data[8,3] = 2967
$data[8,4] = 14313
$data[8,5] = 348
$data[8,6] = 0
$data[8,7] = 488
$data[9,0] = "Turquia"
$data[9,1] = 15679
$data[9,2] = 0
$data[9,3] = 333
$data[9,4] = 15069
$data[9,5] = 847
$data[9,6] = 0
$data[9,7] = 277
$data[10,0] = "Belgica"
$data[10,1] = 13964
$data[10,2] = 0
$data[10,3] = 2132
$data[10,4] = 11004
$data[10,5] = 1088
$data[10,6] = 0
$data[10,7] = 828
$data[11,0] = "Paises Bajos"
$data[11,1] = 13614
$data[11,2] = 0
$data[11,3] = 250
$data[11,4] = 12191
$data[11,5] = 1053
$data[11,6] = 0
$data[11,7] = 1173
$data[12,0] = "Austria"
$data[12,1] = 10711
$data[12,2] = 0
$data[12,3] = 1436
$data[12,4] = 9129
$data[12,5] = 215
$data[12,6] = 0
$data[12,7] = 146
$data[13,0] = "Corea del Sur"
$data[13,1] = 9887
$data[13,2] = 0
$data[13,3] = 5567
$data[13,4] = 4155
$data[13,5] = 55
$data[13,6] = 0
$data[13,7] = 165
$data[14,0] = "Canada"
$data[14,1] = 9731
$data[14,2] = 0
$data[14,3] = 1736
$data[14,4] = 7881
$data[14,5] = 120
$data[14,6] = 0
$data[14,7] = 114
$data[15,0] = "Portugal"
$data[15,1] = 8251
$data[15,2] = 0
$data[15,3] = 43
$data[15,4] = 8021
$data[15,5] = 230
$data[15,6] = 0
$data[15,7] = 187
$data[16,0] = "Brasil"
$data[16,1] = 6931
$data[16,2] = 51
$data[16,3] = 127
$data[16,4] = 6560
$data[16,5] = 296
$data[16,6] = 2
$data[16,7] = 244
$data[17,0] = "Israel"
$data[17,1] = 6092
$data[17,2] = 0
$data[17,3] = 241
$data[17,4] = 5825
$data[17,5] = 95
$data[17,6] = 0
$data[17,7] = 26
$data[18,0] = "Australia"
$data[18,1] = 5105
$data[18,2] = 57
$data[18,3] = 345
$data[18,4] = 4737
$data[18,5] = 50
$data[18,6] = 0
$data[18,7] = 23
$data[19,0] = "Suecia"
$data[19,1] = 4947
$data[19,2] = 0
$data[19,3] = 103
$data[19,4] = 4605
$data[19,5] = 393
$data[19,6] = 0
$data[19,7] = 239
$data[20,0] = "Noruega"
$data[20,1] = 4877
$data[20,2] = 0
$data[20,3] = 13
$data[20,4] = 4820
$data[20,5] = 105
$data[20,6] = 0
$data[20,7] = 44
$data[21,0] = "Chequia"
$data[21,1] = 3589
$data[21,2] = 0
$data[21,3] = 61
$data[21,4] = 3489
$data[21,5] = 70
$data[21,6] = 0
$data[21,7] = 39
$data[22,0] = "Irlanda"
$data[22,1] = 3447
$data[22,2] = 0
$data[22,3] = 5
$data[22,4] = 3357
$data[22,5] = 103
$data[22,6] = 0
$data[22,7] = 85
$data[23,0] = "Dinamarca"
$data[23,1] = 3107
$data[23,2] = 0
$data[23,3] = 894
$data[23,4] = 2109
$data[23,5] = 145
$data[23,6] = 0
$data[23,7] = 104
$data[24,0] = "Chile"
$data[24,1] = 3031
$data[24,2] = 0
$data[24,3] = 234
$data[24,4] = 2781
$data[24,5] = 31
$data[24,6] = 0
$data[24,7] = 16
$data[25,0] = "Malasia"
$data[25,1] = 2908
$data[25,2] = 0
$data[25,3] = 645
$data[25,4] = 2218
$data[25,5] = 102
$data[25,6] = 0
$data[25,7] = 45
$data[26,0] = "Rusia"
$data[26,1] = 2777
$data[26,2] = 0
$data[26,3] = 190
$data[26,4] = 2563
$data[26,5] = 8
$data[26,6] = 0
$data[26,7] = 24
$data[27,0] = "Ecuador"
$data[27,1] = 2758
$data[27,2] = 0
$data[27,3] = 58
$data[27,4] = 2602
$data[27,5] = 100
$data[27,6] = 0
$data[27,7] = 98
$data[28,0] = "Polonia"
$data[28,1] = 2554
$data[28,2] = 0
$data[28,3] = 56
$data[28,4] = 2455
$data[28,5] = 50
$data[28,6] = 0
$data[28,7] = 43
$data[29,0] = "Rumania"
$data[29,1] = 2460
$data[29,2] = 0
$data[29,3] = 252
$data[29,4] = 2116
$data[29,5] = 57
$data[29,6] = 0
$data[29,7] = 92
$data[30,0] = "Japon"
$data[30,1] = 2384
$data[30,2] = 0
$data[30,3] = 472
$data[30,4] = 1855
$data[30,5] = 69
$data[30,6] = 0
$data[30,7] = 57
$data[31,0] = "Luxemburgo"
$data[31,1] = 2319
$data[31,2] = 0
$data[31,3] = 80
$data[31,4] = 2210
$data[31,5] = 31
$data[31,6] = 0
$data[31,7] = 29
$data[32,0] = "Filipinas"
$data[32,1] = 2311
$data[32,2] = 0
$data[32,3] = 50
$data[32,4] = 2165
$data[32,5] = 1
$data[32,6] = 0
$data[32,7] = 96
$data[33,0] = "Pakistan"
$data[33,1] = 2118
$data[33,2] = 0
$data[33,3] = 94
$data[33,4] = 1997
$data[33,5] = 12
$data[33,6] = 0
$data[33,7] = 27
$data[34,0] = "India"
$data[34,1] = 1998
$data[34,2] = 0
$data[34,3] = 148
$data[34,4] = 1792
$data[34,5] = 0
$data[34,6] = 0
$data[34,7] = 58
$data[35,0] = "Tailandia"
$data[35,1] = 1771
$data[35,2] = 0
$data[35,3] = 505
$data[35,4] = 1254
$data[35,5] = 23
$data[35,6] = 0
$data[35,7] = 12
$data[36,0] = "Arabia Saudita"
$data[36,1] = 1720
$data[36,2] = 0
$data[36,3] = 264
$data[36,4] = 1440
$data[36,5] = 31
$data[36,6] = 0
$data[36,7] = 16
$data[37,0] = "Indonesia"
$data[37,1] = 1677
$data[37,2] = 0
$data[37,3] = 103
$data[37,4] = 1417
$data[37,5] = 0
$data[37,6] = 0
$data[37,7] = 157
$data[38,0] = "Finlandia"
$data[38,1] = 1446
$data[38,2] = 0
$data[38,3] = 10
$data[38,4] = 1419
$data[38,5] = 62
$data[38,6] = 0
$data[38,7] = 17
$data[39,0] = "Grecia"
$data[39,1] = 1415
$data[39,2] = 0
$data[39,3] = 52
$data[39,4] = 1312
$data[39,5] = 90
$data[39,6] = 0
$data[39,7] = 51
$data[40,0] = "Sudafrica"
$data[40,1] = 1380
$data[40,2] = 0
$data[40,3] = 50
$data[40,4] = 1325
$data[40,5] = 7
$data[40,6] = 0
$data[40,7] = 5
$data[41,0] = "Peru"
$data[41,1] = 1323
$data[41,2] = 0
$data[41,3] = 394
$data[41,4] = 891
$data[41,5] = 49
$data[41,6] = 0
$data[41,7] = 38
$data[42,0] = "Panama"
$data[42,1] = 1317
$data[42,2] = 0
$data[42,3] = 9
$data[42,4] = 1276
$data[42,5] = 50
$data[42,6] = 0
$data[42,7] = 32
$data[43,0] = "Republica Dominicana"
$data[43,1] = 1284
$data[43,2] = 0
$data[43,3] = 9
$data[43,4] = 1218
$data[43,5] = 0
$data[43,6] = 0
$data[43,7] = 57
$data[44,0] = "Islandia"
$data[44,1] = 1220
$data[44,2] = 0
$data[44,3] = 236
$data[44,4] = 982
$data[44,5] = 12
$data[44,6] = 0
$data[44,7] = 2
$data[45,0] = "Mexico"
$data[45,1] = 1215
$data[45,2] = 0
$data[45,3] = 35
$data[45,4] = 1151
$data[45,5] = 1
$data[45,6] = 0
$data[45,7] = 29
$data[46,0] = "Argentina"
$data[46,1] = 1133
$data[46,2] = 0
$data[46,3] = 248
$data[46,4] = 853
$data[46,5] = 0
$data[46,6] = 0
$data[46,7] = 32
$data[47,0] = "Colombia"
$data[47,1] = 1065
$data[47,2] = 0
$data[47,3] = 39
$data[47,4] = 1009
$data[47,5] = 47
$data[47,6] = 0
$data[47,7] = 17
$data[48,0] = "Serbia"
$data[48,1] = 1060
$data[48,2] = 0
$data[48,3] = 42
$data[48,4] = 990
$data[48,5] = 62
$data[48,6] = 0
$data[48,7] = 28
$data[49,0] = "Singapur"
$data[49,1] = 1000
$data[49,2] = 0
$data[49,3] = 245
$data[49,4] = 752
$data[49,5] = 24
$data[49,6] = 0
$data[49,7] = 3
$data[50,0] = "Croacia"
$data[50,1] = 963
$data[50,2] = 0
$data[50,3] = 73
$data[50,4] = 884
$data[50,5] = 34
$data[50,6] = 0
$data[50,7] = 6
$data[51,0] = "Argelia"
$data[51,1] = 847
$data[51,2] = 0
$data[51,3] = 61
$data[51,4] = 728
$data[51,5] = 0
$data[51,6] = 0
$data[51,7] = 58
$data[52,0] = "Eslovenia"
$data[52,1] = 841
$data[52,2] = 0
$data[52,3] = 10
$data[52,4] = 816
$data[52,5] = 31
$data[52,6] = 0
$data[52,7] = 15
$data[53,0] = "Catar"
$data[53,1] = 835
$data[53,2] = 0
$data[53,3] = 71
$data[53,4] = 762
$data[53,5] = 37
$data[53,6] = 0
$data[53,7] = 2
$data[54,0] = "Emiratos Arabes Unidos"
$data[54,1] = 814
$data[54,2] = 0
$data[54,3] = 61
$data[54,4] = 745
$data[54,5] = 2
$data[54,6] = 0
$data[54,7] = 8
$data[55,0] = "Nueva Zelanda"
$data[55,1] = 797
$data[55,2] = 89
$data[55,3] = 92
$data[55,4] = 704
$data[55,5] = 2
$data[55,6] = 0
$data[55,7] = 1
$data[56,0] = "Ucrania"
$data[56,1] = 794
$data[56,2] = 0
$data[56,3] = 13
$data[56,4] = 761
$data[56,5] = 0
$data[56,6] = 0
$data[56,7] = 20
$data[57,0] = "Estonia"
$data[57,1] = 779
$data[57,2] = 0
$data[57,3] = 33
$data[57,4] = 741
$data[57,5] = 15
$data[57,6] = 0
$data[57,7] = 5
$data[58,0] = "Egipto"
$data[58,1] = 779
$data[58,2] = 0
$data[58,3] = 179
$data[58,4] = 548
$data[58,5] = 0
$data[58,6] = 0
$data[58,7] = 52
$data[59,0] = "Hong Kong"
$data[59,1] = 766
$data[59,2] = 0
$data[59,3] = 147
$data[59,4] = 615
$data[59,5] = 5
$data[59,6] = 0
$data[59,7] = 4
$data[60,0] = "Irak"
$data[60,1] = 728
$data[60,2] = 0
$data[60,3] = 182
$data[60,4] = 494
$data[60,5] = 0
$data[60,6] = 0
$data[60,7] = 52
$data[61,0] = "Crucero"
$data[61,1] = 712
$data[61,2] = 0
$data[61,3] = 603
$data[61,4] = 98
$data[61,5] = 15
$data[61,6] = 0
$data[61,7] = 11
$data[62,0] = "Marruecos"
$data[62,1] = 654
$data[62,2] = 0
$data[62,3] = 29
$data[62,4] = 586
$data[62,5] = 1
$data[62,6] = 0
$data[62,7] = 39
$data[63,0] = "Lituania"
$data[63,1] = 581
$data[63,2] = 0
$data[63,3] = 7
$data[63,4] = 566
$data[63,5] = 11
$data[63,6] = 0
$data[63,7] = 8
$data[64,0] = "Armenia"
$data[64,1] = 571
$data[64,2] = 0
$data[64,3] = 31
$data[64,4] = 536
$data[64,5] = 30
$data[64,6] = 0
$data[64,7] = 4
$data[65,0] = "Barein"
$data[65,1] = 569
$data[65,2] = 0
$data[65,3] = 337
$data[65,4] = 228
$data[65,5] = 3
$data[65,6] = 0
$data[65,7] = 4
$data[66,0] = "Hungria"
$data[66,1] = 525
$data[66,2] = 0
$data[66,3] = 40
$data[66,4] = 465
$data[66,5] = 17
$data[66,6] = 0
$data[66,7] = 20
$data[67,0] = "Libano"
$data[67,1] = 479
$data[67,2] = 0
$data[67,3] = 43
$data[67,4] = 422
$data[67,5] = 5
$data[67,6] = 0
$data[67,7] = 14
$data[68,0] = "Bosnia y Herzegovina"
$data[68,1] = 459
$data[68,2] = 0
$data[68,3] = 19
$data[68,4] = 427
$data[68,5] = 1
$data[68,6] = 0
$data[68,7] = 13
$data[69,0] = "Letonia"
$data[69,1] = 446
$data[69,2] = 0
$data[69,3] = 1
$data[69,4] = 445
$data[69,5] = 3
$data[69,6] = 0
$data[69,7] = 0
$data[70,0] = "Tunez"
$data[70,1] = 423
$data[70,2] = 0
$data[70,3] = 5
$data[70,4] = 406
$data[70,5] = 10
$data[70,6] = 0
$data[70,7] = 12
$data[71,0] = "Moldavia"
$data[71,1] = 423
$data[71,2] = 0
$data[71,3] = 23
$data[71,4] = 395
$data[71,5] = 44
$data[71,6] = 0
$data[71,7] = 5
$data[72,0] = "Bulgaria"
$data[72,1] = 422
$data[72,2] = 0
$data[72,3] = 20
$data[72,4] = 392
$data[72,5] = 18
$data[72,6] = 0
$data[72,7] = 10
$data[73,0] = "Eslovaquia"
$data[73,1] = 400
$data[73,2] = 0
$data[73,3] = 3
$data[73,4] = 396
$data[73,5] = 1
$data[73,6] = 0
$data[73,7] = 1
$data[74,0] = "Principado de Andorra"
$data[74,1] = 390
$data[74,2] = 0
$data[74,3] = 10
$data[74,4] = 366
$data[74,5] = 12
$data[74,6] = 0
$data[74,7] = 14
$data[75,0] = "Kazajistan"
$data[75,1] = 380
$data[75,2] = 0
$data[75,3] = 26
$data[75,4] = 351
$data[75,5] = 6
$data[75,6] = 0
$data[75,7] = 3
$data[76,0] = "Costa Rica"
$data[76,1] = 375
$data[76,2] = 0
$data[76,3] = 4
$data[76,4] = 369
$data[76,5] = 9
$data[76,6] = 0
$data[76,7] = 2
$data[77,0] = "Azerbaiyan"
$data[77,1] = 359
$data[77,2] = 0
$data[77,3] = 26
$data[77,4] = 328
$data[77,5] = 7
$data[77,6] = 0
$data[77,7] = 5
$data[78,0] = "Republica de Macedonia"
$data[78,1] = 354
$data[78,2] = 0
$data[78,3] = 17
$data[78,4] = 326
$data[78,5] = 4
$data[78,6] = 0
$data[78,7] = 11
$data[79,0] = "Uruguay"
$data[79,1] = 350
$data[79,2] = 0
$data[79,3] = 62
$data[79,4] = 286
$data[79,5] = 15
$data[79,6] = 0
$data[79,7] = 2
$data[80,0] = "Taiwan"
$data[80,1] = 329
$data[80,2] = 0
$data[80,3] = 45
$data[80,4] = 279
$data[80,5] = 0
$data[80,6] = 0
$data[80,7] = 5
$data[81,0] = "Republica de Chipre"
$data[81,1] = 320
$data[81,2] = 0
$data[81,3] = 28
$data[81,4] = 283
$data[81,5] = 11
$data[81,6] = 0
$data[81,7] = 9
$data[82,0] = "Kuwait"
$data[82,1] = 317
$data[82,2] = 0
$data[82,3] = 80
$data[82,4] = 237
$data[82,5] = 14
$data[82,6] = 0
$data[82,7] = 0
$data[83,0] = "Burkina Faso"
$data[83,1] = 282
$data[83,2] = 0
$data[83,3] = 46
$data[83,4] = 220
$data[83,5] = 0
$data[83,6] = 0
$data[83,7] = 16
$data[84,0] = "Reunion"
$data[84,1] = 281
$data[84,2] = 0
$data[84,3] = 40
$data[84,4] = 241
$data[84,5] = 3
$data[84,6] = 0
$data[84,7] = 0
$data[85,0] = "Jordania"
$data[85,1] = 278
$data[85,2] = 0
$data[85,3] = 36
$data[85,4] = 237
$data[85,5] = 5
$data[85,6] = 0
$data[85,7] = 5
$data[86,0] = "Albania"
$data[86,1] = 259
$data[86,2] = 0
$data[86,3] = 67
$data[86,4] = 177
$data[86,5] = 7
$data[86,6] = 0
$data[86,7] = 15
$data[87,0] = "Afganistan"
$data[87,1] = 237
$data[87,2] = 0
$data[87,3] = 5
$data[87,4] = 228
$data[87,5] = 0
$data[87,6] = 0
$data[87,7] = 4
$data[88,0] = "San Marino"
$data[88,1] = 236
$data[88,2] = 0
$data[88,3] = 13
$data[88,4] = 195
$data[88,5] = 16
$data[88,6] = 0
$data[88,7] = 28
$data[89,0] = "Camerun"
$data[89,1] = 233
$data[89,2] = 0
$data[89,3] = 10
$data[89,4] = 217
$data[89,5] = 0
$data[89,6] = 0
$data[89,7] = 6
$data[90,0] = "Vietnam"
$data[90,1] = 222
$data[90,2] = 4
$data[90,3] = 63
$data[90,4] = 159
$data[90,5] = 3
$data[90,6] = 0
$data[90,7] = 0
$data[91,0] = "Cuba"
$data[91,1] = 212
$data[91,2] = 0
$data[91,3] = 12
$data[91,4] = 194
$data[91,5] = 3
$data[91,6] = 0
$data[91,7] = 6
$data[92,0] = "Oman"
$data[92,1] = 210
$data[92,2] = 0
$data[92,3] = 34
$data[92,4] = 175
$data[92,5] = 3
$data[92,6] = 0
$data[92,7] = 1
$data[93,0] = "Ghana"
$data[93,1] = 195
$data[93,2] = 0
$data[93,3] = 31
$data[93,4] = 159
$data[93,5] = 1
$data[93,6] = 0
$data[93,7] = 5
$data[94,0] = "Costa de Marfil"
$data[94,1] = 190
$data[94,2] = 0
$data[94,3] = 9
$data[94,4] = 180
$data[94,5] = 0
$data[94,6] = 0
$data[94,7] = 1
$data[95,0] = "Senegal"
$data[95,1] = 190
$data[95,2] = 0
$data[95,3] = 45
$data[95,4] = 144
$data[95,5] = 0
$data[95,6] = 0
$data[95,7] = 1
$data[96,0] = "Malta"
$data[96,1] = 188
$data[96,2] = 0
$data[96,3] = 2
$data[96,4] = 186
$data[96,5] = 2
$data[96,6] = 0
$data[96,7] = 0
$data[97,0] = "Uzbekistan"
$data[97,1] = 181
$data[97,2] = 0
$data[97,3] = 12
$data[97,4] = 167
$data[97,5] = 8
$data[97,6] = 0
$data[97,7] = 2
$data[98,0] = "Nigeria"
$data[98,1] = 174
$data[98,2] = 0
$data[98,3] = 9
$data[98,4] = 163
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 2
$data[99,0] = "Islas Feroe"
$data[99,1] = 173
$data[99,2] = 0
$data[99,3] = 75
$data[99,4] = 98
$data[99,5] = 1
$data[99,6] = 0
$data[99,7] = 0
$data[100,0] = "Honduras"
$data[100,1] = 172
$data[100,2] = 0
$data[100,3] = 3
$data[100,4] = 159
$data[100,5] = 4
$data[100,6] = 0
$data[100,7] = 10
$data[101,0] = "Bielorrusia"
$data[101,1] = 163
$data[101,2] = 0
$data[101,3] = 53
$data[101,4] = 108
$data[101,5] = 2
$data[101,6] = 0
$data[101,7] = 2
$data[102,0] = "Mauricio"
$data[102,1] = 161
$data[102,2] = 0
$data[102,3] = 0
$data[102,4] = 155
$data[102,5] = 1
$data[102,6] = 0
$data[102,7] = 6
$data[103,0] = "Sri Lanka"
$data[103,1] = 146
$data[103,2] = 0
$data[103,3] = 21
$data[103,4] = 122
$data[103,5] = 5
$data[103,6] = 0
$data[103,7] = 3
$data[104,0] = "Venezuela"
$data[104,1] = 144
$data[104,2] = 0
$data[104,3] = 43
$data[104,4] = 98
$data[104,5] = 6
$data[104,6] = 0
$data[104,7] = 3
$data[105,0] = "Martinica"
$data[105,1] = 135
$data[105,2] = 0
$data[105,3] = 27
$data[105,4] = 105
$data[105,5] = 16
$data[105,6] = 0
$data[105,7] = 3
$data[106,0] = "Estado de Palestina"
$data[106,1] = 134
$data[106,2] = 0
$data[106,3] = 18
$data[106,4] = 115
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 1
$data[107,0] = "Brunei"
$data[107,1] = 131
$data[107,2] = 0
$data[107,3] = 52
$data[107,4] = 78
$data[107,5] = 3
$data[107,6] = 0
$data[107,7] = 1
$data[108,0] = "Guadalupe"
$data[108,1] = 125
$data[108,2] = 0
$data[108,3] = 24
$data[108,4] = 95
$data[108,5] = 14
$data[108,6] = 0
$data[108,7] = 6
$data[109,0] = "Montenegro"
$data[109,1] = 123
$data[109,2] = 0
$data[109,3] = 0
$data[109,4] = 121
$data[109,5] = 4
$data[109,6] = 0
$data[109,7] = 2
$data[110,0] = "Georgia"
$data[110,1] = 117
$data[110,2] = 0
$data[110,3] = 23
$data[110,4] = 94
$data[110,5] = 6
$data[110,6] = 0
$data[110,7] = 0
$data[111,0] = "Bolivia"
$data[111,1] = 115
$data[111,2] = 0
$data[111,3] = 1
$data[111,4] = 107
$data[111,5] = 3
$data[111,6] = 0
$data[111,7] = 7
$data[112,0] = "Kirguistan"
$data[112,1] = 111
$data[112,2] = 0
$data[112,3] = 3
$data[112,4] = 108
$data[112,5] = 3
$data[112,6] = 0
$data[112,7] = 0
$data[113,0] = "Consejo Danes para los Refugiados"
$data[113,1] = 109
$data[113,2] = 0
$data[113,3] = 3
$data[113,4] = 97
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 9
$data[114,0] = "Camboya"
$data[114,1] = 109
$data[114,2] = 0
$data[114,3] = 25
$data[114,4] = 84
$data[114,5] = 1
$data[114,6] = 0
$data[114,7] = 0
$data[115,0] = "Mayotte"
$data[115,1] = 101
$data[115,2] = 0
$data[115,3] = 10
$data[115,4] = 90
$data[115,5] = 3
$data[115,6] = 0
$data[115,7] = 1
$data[116,0] = "Trinidad yTobago"
$data[116,1] = 90
$data[116,2] = 0
$data[116,3] = 1
$data[116,4] = 84
$data[116,5] = 0
$data[116,6] = 0
$data[116,7] = 5
$data[117,0] = "Ruanda"
$data[117,1] = 82
$data[117,2] = 0
$data[117,3] = 0
$data[117,4] = 82
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 0
$data[118,0] = "Kenia"
$data[118,1] = 81
$data[118,2] = 0
$data[118,3] = 3
$data[118,4] = 77
$data[118,5] = 2
$data[118,6] = 0
$data[118,7] = 1
$data[119,0] = "Gibraltar"
$data[119,1] = 81
$data[119,2] = 0
$data[119,3] = 34
$data[119,4] = 47
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 0
$data[120,0] = "Paraguay"
$data[120,1] = 77
$data[120,2] = 8
$data[120,3] = 2
$data[120,4] = 72
$data[120,5] = 4
$data[120,6] = 0
$data[120,7] = 3
$data[121,0] = "Niger"
$data[121,1] = 74
$data[121,2] = 0
$data[121,3] = 0
$data[121,4] = 69
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 5
$data[122,0] = "Liechtenstein"
$data[122,1] = 72
$data[122,2] = 0
$data[122,3] = 0
$data[122,4] = 72
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 0
$data[123,0] = "Isla de Man"
$data[123,1] = 68
$data[123,2] = 0
$data[123,3] = 0
$data[123,4] = 67
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 1
$data[124,0] = "Madagascar"
$data[124,1] = 57
$data[124,2] = 0
$data[124,3] = 0
$data[124,4] = 57
$data[124,5] = 6
$data[124,6] = 0
$data[124,7] = 0
$data[125,0] = "Aruba"
$data[125,1] = 55
$data[125,2] = 0
$data[125,3] = 1
$data[125,4] = 54
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 0
$data[126,0] = "Monaco"
$data[126,1] = 55
$data[126,2] = 0
$data[126,3] = 2
$data[126,4] = 52
$data[126,5] = 2
$data[126,6] = 0
$data[126,7] = 1
$data[127,0] = "Banglades"
$data[127,1] = 54
$data[127,2] = 0
$data[127,3] = 25
$data[127,4] = 23
$data[127,5] = 1
$data[127,6] = 0
$data[127,7] = 6
$data[128,0] = "Guayana Francesa"
$data[128,1] = 51
$data[128,2] = 0
$data[128,3] = 15
$data[128,4] = 36
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 0
$data[129,0] = "Guatemala"
$data[129,1] = 46
$data[129,2] = 7
$data[129,3] = 12
$data[129,4] = 33
$data[129,5] = 1
$data[129,6] = 0
$data[129,7] = 1
$data[130,0] = "Barbados"
$data[130,1] = 45
$data[130,2] = 0
$data[130,3] = 0
$data[130,4] = 45
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 0
$data[131,0] = "Uganda"
$data[131,1] = 44
$data[131,2] = 0
$data[131,3] = 0
$data[131,4] = 44
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 0
$data[132,0] = "Jamaica"
$data[132,1] = 44
$data[132,2] = 0
$data[132,3] = 2
$data[132,4] = 39
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 3
$data[133,0] = "Macao"
$data[133,1] = 41
$data[133,2] = 0
$data[133,3] = 10
$data[133,4] = 31
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 0
$data[134,0] = "Puerto Rico"
$data[134,1] = 39
$data[134,2] = 0
$data[134,3] = 1
$data[134,4] = 36
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 2
$data[135,0] = "Polinesia Francesa"
$data[135,1] = 37
$data[135,2] = 0
$data[135,3] = 0
$data[135,4] = 37
$data[135,5] = 1
$data[135,6] = 0
$data[135,7] = 0
$data[136,0] = "Zambia"
$data[136,1] = 36
$data[136,2] = 0
$data[136,3] = 0
$data[136,4] = 36
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 0
$data[137,0] = "Togo"
$data[137,1] = 36
$data[137,2] = 0
$data[137,3] = 10
$data[137,4] = 24
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 2
$data[138,0] = "Republica de Yibuti"
$data[138,1] = 33
$data[138,2] = 0
$data[138,3] = 0
$data[138,4] = 33
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 0
$data[139,0] = "El Salvador"
$data[139,1] = 33
$data[139,2] = 0
$data[139,3] = 0
$data[139,4] = 31
$data[139,5] = 4
$data[139,6] = 0
$data[139,7] = 2
$data[140,0] = "Guam"
$data[140,1] = 32
$data[140,2] = 0
$data[140,3] = 0
$data[140,4] = 31
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 1
$data[141,0] = "Bermudas"
$data[141,1] = 32
$data[141,2] = 0
$data[141,3] = 10
$data[141,4] = 22
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 0
$data[142,0] = "Mali"
$data[142,1] = 31
$data[142,2] = 0
$data[142,3] = 0
$data[142,4] = 28
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 3
$data[143,0] = "Guinea"
$data[143,1] = 30
$data[143,2] = 0
$data[143,3] = 0
$data[143,4] = 30
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 0
$data[144,0] = "Etiopia"
$data[144,1] = 29
$data[144,2] = 0
$data[144,3] = 2
$data[144,4] = 27
$data[144,5] = 2
$data[144,6] = 0
$data[144,7] = 0
$data[145,0] = "Islas Caimanes"
$data[145,1] = 22
$data[145,2] = 0
$data[145,3] = 0
$data[145,4] = 21
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 1
$data[146,0] = "Congo"
$data[146,1] = 22
$data[146,2] = 0
$data[146,3] = 0
$data[146,4] = 20
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 2
$data[147,0] = "San Martin (Parte Francesa)"
$data[147,1] = 22
$data[147,2] = 0
$data[147,3] = 2
$data[147,4] = 19
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 1
$data[148,0] = "Bahamas"
$data[148,1] = 21
$data[148,2] = 0
$data[148,3] = 1
$data[148,4] = 19
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 1
$data[149,0] = "Tanzania"
$data[149,1] = 20
$data[149,2] = 0
$data[149,3] = 1
$data[149,4] = 18
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 1
$data[150,0] = "Guyana"
$data[150,1] = 19
$data[150,2] = 0
$data[150,3] = 0
$data[150,4] = 16
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 3
$data[151,0] = "Maldivas"
$data[151,1] = 19
$data[151,2] = 0
$data[151,3] = 13
$data[151,4] = 6
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 0
$data[152,0] = "Gabon"
$data[152,1] = 18
$data[152,2] = 0
$data[152,3] = 0
$data[152,4] = 17
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 1
$data[153,0] = "Islas Virgenes de los Estados Unidos"
$data[153,1] = 17
$data[153,2] = 0
$data[153,3] = 0
$data[153,4] = 17
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 0
$data[154,0] = "Birmania"
$data[154,1] = 16
$data[154,2] = 0
$data[154,3] = 0
$data[154,4] = 15
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 1
$data[155,0] = "Nueva Caledonia"
$data[155,1] = 16
$data[155,2] = 0
$data[155,3] = 1
$data[155,4] = 15
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 0
$data[156,0] = "Haiti"
$data[156,1] = 16
$data[156,2] = 0
$data[156,3] = 1
$data[156,4] = 15
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 0
$data[157,0] = "San Martin (Parte Holandesa)"
$data[157,1] = 16
$data[157,2] = 0
$data[157,3] = 6
$data[157,4] = 9
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 1
$data[158,0] = "Eritrea"
$data[158,1] = 15
$data[158,2] = 0
$data[158,3] = 0
$data[158,4] = 15
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 0
$data[159,0] = "Guinea Ecuatorial"
$data[159,1] = 15
$data[159,2] = 0
$data[159,3] = 1
$data[159,4] = 14
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 0
$data[160,0] = "Mongolia"
$data[160,1] = 14
$data[160,2] = 0
$data[160,3] = 2
$data[160,4] = 12
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 0
$data[161,0] = "Namibia"
$data[161,1] = 14
$data[161,2] = 0
$data[161,3] = 2
$data[161,4] = 12
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 0
$data[162,0] = "Santa Lucia"
$data[162,1] = 13
$data[162,2] = 0
$data[162,3] = 1
$data[162,4] = 12
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 0
$data[163,0] = "Benin"
$data[163,1] = 13
$data[163,2] = 0
$data[163,3] = 1
$data[163,4] = 12
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 0
$data[164,0] = "Dominica"
$data[164,1] = 12
$data[164,2] = 0
$data[164,3] = 0
$data[164,4] = 12
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 0
$data[165,0] = "Curazao"
$data[165,1] = 11
$data[165,2] = 0
$data[165,3] = 3
$data[165,4] = 7
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 1
$data[166,0] = "Mozambique"
$data[166,1] = 10
$data[166,2] = 0
$data[166,3] = 0
$data[166,4] = 10
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 0
$data[167,0] = "Seychelles"
$data[167,1] = 10
$data[167,2] = 0
$data[167,3] = 0
$data[167,4] = 10
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 0
$data[168,0] = "Libia"
$data[168,1] = 10
$data[168,2] = 0
$data[168,3] = 0
$data[168,4] = 10
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 0
$data[169,0] = "Granada"
$data[169,1] = 10
$data[169,2] = 1
$data[169,3] = 0
$data[169,4] = 10
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 0
$data[170,0] = "Laos"
$data[170,1] = 10
$data[170,2] = 0
$data[170,3] = 0
$data[170,4] = 10
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 0
$data[171,0] = "Surinam"
$data[171,1] = 10
$data[171,2] = 0
$data[171,3] = 0
$data[171,4] = 10
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 0
$data[172,0] = "Groenlandia"
$data[172,1] = 10
$data[172,2] = 0
$data[172,3] = 2
$data[172,4] = 8
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 0
$data[173,0] = "Siria"
$data[173,1] = 10
$data[173,2] = 0
$data[173,3] = 0
$data[173,4] = 8
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 2
$data[174,0] = "Guinea-Bisau"
$data[174,1] = 9
$data[174,2] = 0
$data[174,3] = 0
$data[174,4] = 9
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 0
$data[175,0] = "Suazilandia"
$data[175,1] = 9
$data[175,2] = 0
$data[175,3] = 0
$data[175,4] = 9
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = "Montserrat"
$data[176,1] = 9
$data[176,2] = 0
$data[176,3] = 0
$data[176,4] = 7
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 2
$data[177,0] = "San Cristobal y Nieves"
$data[177,1] = 8
$data[177,2] = 0
$data[177,3] = 0
$data[177,4] = 8
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = "Zimbabue"
$data[178,1] = 8
$data[178,2] = 0
$data[178,3] = 0
$data[178,4] = 7
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 1
$data[179,0] = "Angola"
$data[179,1] = 8
$data[179,2] = 0
$data[179,3] = 1
$data[179,4] = 5
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 2
$data[180,0] = "Antigua y Barbuda"
$data[180,1] = 7
$data[180,2] = 0
$data[180,3] = 0
$data[180,4] = 7
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = "Republica del Chad"
$data[181,1] = 7
$data[181,2] = 0
$data[181,3] = 0
$data[181,4] = 7
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 0
$data[182,0] = "Sudan"
$data[182,1] = 7
$data[182,2] = 0
$data[182,3] = 2
$data[182,4] = 3
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 2
$data[183,0] = "Liberia"
$data[183,1] = 6
$data[183,2] = 0
$data[183,3] = 0
$data[183,4] = 6
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = "Islas Turcas y Caicos"
$data[184,1] = 6
$data[184,2] = 0
$data[184,3] = 0
$data[184,4] = 6
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = "Santa Sede"
$data[185,1] = 6
$data[185,2] = 0
$data[185,3] = 0
$data[185,4] = 6
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = "Cabo Verde"
$data[186,1] = 6
$data[186,2] = 0
$data[186,3] = 0
$data[186,4] = 5
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 1
$data[187,0] = "San Bartolome"
$data[187,1] = 6
$data[187,2] = 0
$data[187,3] = 1
$data[187,4] = 5
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 0
$data[188,0] = "Mauritania"
$data[188,1] = 6
$data[188,2] = 0
$data[188,3] = 2
$data[188,4] = 3
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 1
$data[189,0] = "Fiyi"
$data[189,1] = 5
$data[189,2] = 0
$data[189,3] = 0
$data[189,4] = 5
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 0
$data[190,0] = "Nicaragua"
$data[190,1] = 5
$data[190,2] = 0
$data[190,3] = 0
$data[190,4] = 4
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 1
$data[191,0] = "Somalia"
$data[191,1] = 5
$data[191,2] = 0
$data[191,3] = 1
$data[191,4] = 4
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 0
$data[192,0] = "Nepal"
$data[192,1] = 5
$data[192,2] = 0
$data[192,3] = 1
$data[192,4] = 4
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = "Butan"
$data[193,1] = 4
$data[193,2] = 0
$data[193,3] = 0
$data[193,4] = 4
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 0
$data[194,0] = "Botsuana"
$data[194,1] = 4
$data[194,2] = 0
$data[194,3] = 0
$data[194,4] = 3
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 1
$data[195,0] = "Gambia"
$data[195,1] = 4
$data[195,2] = 0
$data[195,3] = 2
$data[195,4] = 1
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 1
$data[196,0] = "Belice"
$data[196,1] = 3
$data[196,2] = 0
$data[196,3] = 0
$data[196,4] = 3
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 0
$data[197,0] = "Republica de Africa Central"
$data[197,1] = 3
$data[197,2] = 0
$data[197,3] = 0
$data[197,4] = 3
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0
$data[198,0] = "Islas Virgenes Britanicas"
$data[198,1] = 3
$data[198,2] = 0
$data[198,3] = 0
$data[198,4] = 3
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0
$data[199,0] = "Bonaire, San Eustaquio y Saba"
$data[199,1] = 2
$data[199,2] = 0
$data[199,3] = 0
$data[199,4] = 2
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 0
$data[200,0] = "Anguila"
$data[200,1] = 2
$data[200,2] = 0
$data[200,3] = 0
$data[200,4] = 2
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0
$data[201,0] = "Sierra Leona"
$data[201,1] = 2
$data[201,2] = 0
$data[201,3] = 0
$data[201,4] = 2
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0
$data[202,0] = "Burundi"
$data[202,1] = 2
$data[202,2] = 0
$data[202,3] = 0
$data[202,4] = 2
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 0
$data[203,0] = "San Vicente y las Granadinas"
$data[203,1] = 2
$data[203,2] = 0
$data[203,3] = 1
$data[203,4] = 1
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 0
$data[204,0] = "Timor Oriental"
$data[204,1] = 1
$data[204,2] = 0
$data[204,3] = 0
$data[204,4] = 1
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 0
$data[205,0] = "Papua Nueva Guinea"
$data[205,1] = 1
$data[205,2] = 0
$data[205,3] = 0
$data[205,4] = 1
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0

$ws.Range("A4:H209").Value = $data

Write-Output "Updated 206 country rows (A4:H209) and refreshed A1 timestamp."
